$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1 and Q1 with the same style as O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: swap I/K and M/O, then add P and Q columns with value 2
for ($r = 2; $r -le 25; $r++) {
    $i = $ws.Cells.Item($r, 9).Value()
    $k = $ws.Cells.Item($r, 11).Value()
    $ws.Cells.Item($r, 9).Value = $k
    $ws.Cells.Item($r, 11).Value = $i

    $m = $ws.Cells.Item($r, 13).Value()
    $o = $ws.Cells.Item($r, 15).Value()
    $ws.Cells.Item($r, 13).Value = $o
    $ws.Cells.Item($r, 15).Value = $m

    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
